$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 253
$ws1.Range("F5").Value = 2979
$ws1.Range("F6").Value = 2011
$ws1.Range("F9").Value = 1107
$ws1.Range("F11").Value = 642

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 253
$ws4.Range("F5").Value = 2979
$ws4.Range("F6").Value = 2011
$ws4.Range("F10").Value = 1107
$ws4.Range("F12").Value = 642

$wb.Save()
